$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New daily-work-progress rows (26-28), modelled on the existing rows so the
# date/time/border/alignment formatting matches the rest of the table.
# ---------------------------------------------------------------------------

# Row 26 - reuse row 25's per-column formatting, then overwrite the values.
$ws.Range("A25").Copy() | Out-Null
$ws.Range("A26").PasteSpecial(-4122) | Out-Null
$ws.Range("B25").Copy() | Out-Null
$ws.Range("B26").PasteSpecial(-4122) | Out-Null
$ws.Range("C25:D25").Copy() | Out-Null
$ws.Range("C26:D26").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null

$ws.Range("A26").Value = 44945
$ws.Range("B26").Value = "Regulator Structural Design Day-5 Coding"
$ws.Range("C26").Value = 0.39583333333333331
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = "A.K.M Saifuddin"

# Row 27 - same template.
$ws.Range("A25").Copy() | Out-Null
$ws.Range("A27").PasteSpecial(-4122) | Out-Null
$ws.Range("B25").Copy() | Out-Null
$ws.Range("B27").PasteSpecial(-4122) | Out-Null
$ws.Range("C25:D25").Copy() | Out-Null
$ws.Range("C27:D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null

$ws.Range("A27").Value = 44946
$ws.Range("B27").Value = "Regulator Structural Design Day-6 Coding"
$ws.Range("C27").Value = 0.39583333333333331
$ws.Range("D27").Value = 0.99930555555555556
$ws.Range("E27").Value = "A.K.M Saifuddin"

# Row 28 - A/B/E follow the same template; C28 gets a brand new (borderless)
# time style and D28 is left blank.
$ws.Range("A25").Copy() | Out-Null
$ws.Range("A28").PasteSpecial(-4122) | Out-Null
$ws.Range("B25").Copy() | Out-Null
$ws.Range("B28").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null

$ws.Range("A28").Value = 44947
$ws.Range("B28").Value = "Regulator Structural Design Day-7 Coding"
$ws.Range("E28").Value = "A.K.M Saifuddin"

$ws.Range("C28").Value = 0
$ws.Range("C28").NumberFormat = "h:mm"
$ws.Range("C28").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Update the saved view: drop the old scrolled/selected cell and select B32.
# ---------------------------------------------------------------------------
$ws.Range("B32").Select() | Out-Null
